# Insert a new weekly price-report row for "Haba" at row 98 (pushing the
# existing rows 98-111 down to 99-112), then populate the new row with the
# latest observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(98).Insert()

$ws.Cells.Item(98, 1).Value = 6
$ws.Cells.Item(98, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(98, 3).Value = "Metropolitana"
$ws.Cells.Item(98, 4).Value = 44449
$ws.Cells.Item(98, 5).Value = 13
$ws.Cells.Item(98, 6).Value = 100112026
$ws.Cells.Item(98, 7).Value = "Haba"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 200
$ws.Cells.Item(98, 11).Value = 14000
$ws.Cells.Item(98, 12).Value = 15000
$ws.Cells.Item(98, 13).Value = 14400
$ws.Cells.Item(98, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(98, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(98, 16).Value = 576
$ws.Cells.Item(98, 17).Value = 25
$ws.Cells.Item(98, 18).Value = "Hortaliza"
